# Auto-generated: apply scheduled market-price / profit refresh to Seraph_Profits sheets
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per-row across all 8 job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2936.7778
$ws.Range("I17").Value = 1300
$ws.Range("J17").Value = 3033.0588
$ws.Range("K17").Value = 3900
$ws.Range("L17").Value = 9099.1764
$ws.Range("M17").Value = -3732
$ws.Range("N17").Value = -9435.1764
$ws.Range("H33").Value = 196
$ws.Range("I33").Value = 70.833336
$ws.Range("K33").Value = 70.833336
$ws.Range("M33").Value = 158.166664
$ws.Range("H58").Value = 3394.2856
$ws.Range("J58").Value = 5900
$ws.Range("L58").Value = 17700
$ws.Range("N58").Value = -18000
$ws.Range("H64").Value = 4723.75
$ws.Range("I64").Value = 4396.6665
$ws.Range("K64").Value = 4396.6665
$ws.Range("M64").Value = -4148.6665
$ws.Range("H67").Value = 4723.75
$ws.Range("I67").Value = 4396.6665
$ws.Range("K67").Value = 4396.6665
$ws.Range("M67").Value = -3538.6665
$ws.Range("H74").Value = 2999.5
$ws.Range("I74").Value = 2999.5
$ws.Range("K74").Value = 2999.5
$ws.Range("M74").Value = -2063.5
$ws.Range("H77").Value = 2999.5
$ws.Range("I77").Value = 2999.5
$ws.Range("K77").Value = 14997.5
$ws.Range("M77").Value = -10317.5
$ws.Range("H80").Value = 14557.077
$ws.Range("I80").Value = 378.57144
$ws.Range("J80").Value = 31098.666
$ws.Range("K80").Value = 1135.71432
$ws.Range("L80").Value = 93295.99800000001
$ws.Range("M80").Value = -137.71432
$ws.Range("N80").Value = -95291.99800000001
$ws.Range("H83").Value = 14557.077
$ws.Range("I83").Value = 378.57144
$ws.Range("J83").Value = 31098.666
$ws.Range("K83").Value = 3407.14296
$ws.Range("L83").Value = 279887.994
$ws.Range("M83").Value = 1584.85704
$ws.Range("N83").Value = -289871.994
$ws.Range("H88").Value = 1580.9231
$ws.Range("J88").Value = 2483.2856
$ws.Range("L88").Value = 2483.2856
$ws.Range("N88").Value = -3295.2856
$ws.Range("H91").Value = 1580.9231
$ws.Range("J91").Value = 2483.2856
$ws.Range("L91").Value = 2483.2856
$ws.Range("N91").Value = -5291.2856
$ws.Range("H112").Value = 2302.5
$ws.Range("J112").Value = 2586.6316
$ws.Range("L112").Value = 7759.8948
$ws.Range("N112").Value = -9975.8948
$ws.Range("H113").Value = 3040.7646
$ws.Range("J113").Value = 3524.75
$ws.Range("L113").Value = 3524.75
$ws.Range("N113").Value = -10032.75
$ws.Range("H129").Value = 3192
$ws.Range("I129").Value = 1568.5714
$ws.Range("J129").Value = 5086
$ws.Range("K129").Value = 4705.7142
$ws.Range("L129").Value = 15258
$ws.Range("M129").Value = 294.2857999999997
$ws.Range("N129").Value = -25258
$ws.Range("H131").Value = 7911.75
$ws.Range("I131").Value = 1647
$ws.Range("K131").Value = 4941
$ws.Range("M131").Value = 99
$ws.Range("H136").Value = 95000
$ws.Range("J136").Value = 95000
$ws.Range("L136").Value = 95000
$ws.Range("N136").Value = -105200
$ws.Range("H137").Value = 3417.6667
$ws.Range("I137").Value = 2264.5
$ws.Range("J137").Value = 3878.9333
$ws.Range("K137").Value = 6793.5
$ws.Range("L137").Value = 11636.7999
$ws.Range("M137").Value = -4243.5
$ws.Range("N137").Value = -16736.7999
$ws.Range("H138").Value = 16118.9
$ws.Range("I138").Value = 13098.2
$ws.Range("K138").Value = 39294.60000000001
$ws.Range("M138").Value = -34154.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6586.8237
$ws.Range("I32").Value = 2740.037
$ws.Range("J32").Value = 21424.428
$ws.Range("K32").Value = 2740.037
$ws.Range("L32").Value = 21424.428
$ws.Range("M32").Value = -2453.037
$ws.Range("N32").Value = -21998.428
$ws.Range("H45").Value = 1946.6
$ws.Range("I45").Value = 1911
$ws.Range("K45").Value = 1911
$ws.Range("M45").Value = -1534
$ws.Range("H61").Value = 1520.1538
$ws.Range("I61").Value = 1520.1538
$ws.Range("K61").Value = 1520.1538
$ws.Range("M61").Value = -1308.1538
$ws.Range("H122").Value = 335952.72
$ws.Range("J122").Value = 4530.4443
$ws.Range("L122").Value = 13591.3329
$ws.Range("N122").Value = -18491.3329
$ws.Range("H136").Value = 1520.1538
$ws.Range("I136").Value = 1520.1538
$ws.Range("K136").Value = 4560.4614
$ws.Range("M136").Value = -2010.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2206.6
$ws.Range("I107").Value = 1595.75
$ws.Range("K107").Value = 1595.75
$ws.Range("M107").Value = 324.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3826.4348
$ws.Range("I31").Value = 1913.875
$ws.Range("J31").Value = 4846.467
$ws.Range("K31").Value = 1913.875
$ws.Range("L31").Value = 4846.467
$ws.Range("M31").Value = -1618.875
$ws.Range("N31").Value = -5436.467
$ws.Range("H34").Value = 3826.4348
$ws.Range("I34").Value = 1913.875
$ws.Range("J34").Value = 4846.467
$ws.Range("K34").Value = 1913.875
$ws.Range("L34").Value = 4846.467
$ws.Range("M34").Value = -1711.875
$ws.Range("N34").Value = -5250.467
$ws.Range("H132").Value = 3602.9656
$ws.Range("I132").Value = 2988.318
$ws.Range("K132").Value = 8964.954000000002
$ws.Range("M132").Value = -6434.954000000002
$ws.Range("H134").Value = 4443.6
$ws.Range("I134").Value = 3204
$ws.Range("K134").Value = 9612
$ws.Range("M134").Value = -7077

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 7469
$ws.Range("I59").Value = 7469
$ws.Range("K59").Value = 22407
$ws.Range("M59").Value = -21867
$ws.Range("H86").Value = 456.66666
$ws.Range("J86").Value = 487.5
$ws.Range("L86").Value = 1462.5
$ws.Range("N86").Value = -3834.5
$ws.Range("H89").Value = 456.66666
$ws.Range("J89").Value = 487.5
$ws.Range("L89").Value = 4387.5
$ws.Range("N89").Value = -16243.5
$ws.Range("H114").Value = 665.3
$ws.Range("I114").Value = 412.83334
$ws.Range("J114").Value = 1044
$ws.Range("K114").Value = 1238.50002
$ws.Range("L114").Value = 3132
$ws.Range("M114").Value = 2015.49998
$ws.Range("N114").Value = -9640
$ws.Range("H131").Value = 1650.9231
$ws.Range("I131").Value = 641.3333
$ws.Range("J131").Value = 1953.8
$ws.Range("K131").Value = 1923.9999
$ws.Range("L131").Value = 5861.4
$ws.Range("M131").Value = 3116.0001
$ws.Range("N131").Value = -15941.4
$ws.Range("H140").Value = 2381.0667
$ws.Range("I140").Value = 2381.0667
$ws.Range("K140").Value = 7143.2001
$ws.Range("M140").Value = -1963.2001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3168.2856
$ws.Range("I80").Value = 2872.6667
$ws.Range("J80").Value = 3390
$ws.Range("K80").Value = 2872.6667
$ws.Range("L80").Value = 3390
$ws.Range("M80").Value = -1874.6667
$ws.Range("N80").Value = -5386
$ws.Range("H83").Value = 3168.2856
$ws.Range("I83").Value = 2872.6667
$ws.Range("J83").Value = 3390
$ws.Range("K83").Value = 14363.3335
$ws.Range("L83").Value = 16950
$ws.Range("M83").Value = -9371.333500000001
$ws.Range("N83").Value = -26934
$ws.Range("H122").Value = 919963.9399999999
$ws.Range("J122").Value = 2004796.6
$ws.Range("L122").Value = 6014389.800000001
$ws.Range("N122").Value = -6019289.800000001
$ws.Range("H132").Value = 4266.3335
$ws.Range("I132").Value = 3918.5
$ws.Range("K132").Value = 11755.5
$ws.Range("M132").Value = -9225.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 967.8182
$ws.Range("I93").Value = 785.7143
$ws.Range("J93").Value = 1286.5
$ws.Range("K93").Value = 785.7143
$ws.Range("L93").Value = 1286.5
$ws.Range("M93").Value = 462.2857
$ws.Range("N93").Value = -3782.5
$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524
$ws.Range("H122").Value = 4266.5
$ws.Range("I122").Value = 3919.8
$ws.Range("K122").Value = 11759.4
$ws.Range("M122").Value = -9309.400000000001
$ws.Range("H136").Value = 4160.3335
$ws.Range("J136").Value = 5000
$ws.Range("L136").Value = 15000
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2200.8635
$ws.Range("I81").Value = 1920.95
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 3841.9
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -2780.9
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 2200.8635
$ws.Range("I84").Value = 1920.95
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 19209.5
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -13905.5
$ws.Range("N84").Value = -60608
$ws.Range("H132").Value = 1181.7084
$ws.Range("I132").Value = 718.15
$ws.Range("K132").Value = 2154.45
$ws.Range("M132").Value = 375.5500000000002
